$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 37; this pushes the existing rows 37-70
# down to rows 39-72 (dimension becomes A1:T72) while keeping their
# formatting (style s="2" on column D carries over to the new rows too).
$ws.Range("A37:A38").EntireRow.Insert()

# --- New row 37 ---
$ws.Range("A37").Value = 10
$ws.Range("B37").Value = "Vega Modelo de Temuco"
$ws.Range("C37").Value = "La Araucanía"
$ws.Range("D37").Value = (Get-Date -Year 2021 -Month 12 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E37").Value = 9
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100101
$ws.Range("H37").Value = "Berries"
$ws.Range("I37").Value = 100101001
$ws.Range("J37").Value = "Arándano (blue)"
$ws.Range("K37").Value = "Sin especificar"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 350
$ws.Range("N37").Value = 2800
$ws.Range("O37").Value = 3000
$ws.Range("P37").Value = 2886
$ws.Range("Q37").Value = "$/kilo"
$ws.Range("R37").Value = "Región del Maule"
$ws.Range("S37").Value = 2886
$ws.Range("T37").Value = 1

# --- New row 38 ---
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = (Get-Date -Year 2021 -Month 12 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100101
$ws.Range("H38").Value = "Berries"
$ws.Range("I38").Value = 100101001
$ws.Range("J38").Value = "Arándano (blue)"
$ws.Range("K38").Value = "Sin especificar"
$ws.Range("L38").Value = "Segunda"
$ws.Range("M38").Value = 300
$ws.Range("N38").Value = 2300
$ws.Range("O38").Value = 2300
$ws.Range("P38").Value = 2300
$ws.Range("Q38").Value = "$/kilo"
$ws.Range("R38").Value = "Región del Maule"
$ws.Range("S38").Value = 2300
$ws.Range("T38").Value = 1
